$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mojibake text in column B (Tweet) - doubled 'EE' sequences
$ws.Cells.Item(33, 2).Value = "| At home sick... ðŸŽ¼The bluesðŸŽ¼ won't cure it so I need ideas ðŸŽ¸ðŸˁEEEE­ | #sorethroat #sick #blues #music #fallweather #carletonuniversity #ottawa"
$ws.Cells.Item(38, 2).Value = "Thought I had a pretty solid GPA as a kin major and now that I look at the average for dpt programs I feel even more discouraged ðŸˁEEEEª"
$ws.Cells.Item(39, 2).Value = "depress ðŸˁEEEE"
$ws.Cells.Item(49, 2).Value = "Am I the only person who dislikes fall? #FirstDayofFall #leaves #thingsdie #depressing #cold #noflipflops ðŸ�EEEE��EEEEŽðŸEEEE¾ðŸ�EEEE��EEEEŽðŸEEEE½ðŸ�EEEE��EEEEŽðŸEEEE»ðŸEEEE�EEEE��EEEEðŸ"
$ws.Cells.Item(62, 2).Value = "If anybody needs me I'll be drowning my blues in a sea of whiskey ðŸEEEE»"
$ws.Cells.Item(63, 2).Value = "@Gaychel22 @Kellee_II @seattlepi even hard facts don't seem to be sinking in... I despair ðŸˁEEEE©"
$ws.Cells.Item(65, 2).Value = "It's just begun ladies and gents. The war on racism has begun. #sadness #Segregation continues.... smh ðŸˁEEEE�EEEE��EEEE"
$ws.Cells.Item(67, 2).Value = "Baaarissshhhhh + sad song =  prefect night â�EEEE�EEEE��EEEE��EEEE feeling alone"
$ws.Cells.Item(88, 2).Value = "One step forward, two steps backward, the link to RogerFedererShop doesnÁEEEE´t work.ðŸˁEEEE° I am losing hope about Roger Federer new Website #sadness"
$ws.Cells.Item(116, 2).Value = "Wow the #Denver housing market is #depressing no places to buy in my range ðŸˁEEEE¢ #FML"
$ws.Cells.Item(138, 2).Value = "Can I just sulk in peace ðŸˁEEEE�EEEE��EEEE"
$ws.Cells.Item(157, 2).Value = "@lucy_hyner @Soulboy2266 sadly not !! One less hour drinking time ðŸˁEEEE¢ðŸEEEE»"
$ws.Cells.Item(180, 2).Value = "Even a pencilâŁEEEEEEEE never #stayed  with me until it's #end âš« ðŸˁEEEEž"
$ws.Cells.Item(186, 2).Value = "The fact I haven't had to wear a bra for a week and knowing I'll have to start wearing one again after tomorrow is depressing ðŸ�EEEE�EEEE��EEEE��EEEEðŸ�EEEE�EEEE��EEEE��EEEE"
$ws.Cells.Item(187, 2).Value = "Liam is too distant makes me mourn ðŸˁEEEEª"
$ws.Cells.Item(249, 2).Value = "@iTriborg â�EEEE�EEEE��EEEE��EEEE make him feel vigorous. 'Fine. You can kill me now.' Said Hestia with a display of only despair rather than her joyful â�EEEE�EEEE��EEEE��EEEE"
$ws.Cells.Item(261, 2).Value = "@JohnWildy71 something, too confused to type an entire word ðŸˁEEEEŁEEEE See, I can laugh again. My hour of sadness has almost passed"
$ws.Cells.Item(265, 2).Value = "It is a solemn thing, and no small scandal in the Kingdom, to see Godâ�EEEE�EEEE��EEEE�EEEE�s children starving while seated at the Fatherâ�EEEE�EEEE��EEEE�EEEE�s table. -AW Tozer"
$ws.Cells.Item(272, 2).Value = "Why is it that we rejoice at a birth and grieve at a funeral? It is because we are not the person involved. â�EEEE�EEEE��EEEE�� Mark Twain"

# Trim trailing-space duplicate annotation strings in column F (anno2_e)
$ws.Cells.Item(120, 6).Value = "displeased"
$ws.Cells.Item(125, 6).Value = "displeased"
$ws.Cells.Item(133, 6).Value = "displeased"
$ws.Cells.Item(135, 6).Value = "depressed"
$ws.Cells.Item(140, 6).Value = "displeased"
$ws.Cells.Item(141, 6).Value = "depressed"
$ws.Cells.Item(157, 6).Value = "depressed"
$ws.Cells.Item(161, 6).Value = "depressed"
$ws.Cells.Item(172, 6).Value = "depressed"
$ws.Cells.Item(173, 6).Value = "depressed"
$ws.Cells.Item(191, 6).Value = "depressed"
$ws.Cells.Item(195, 6).Value = "depressed"
$ws.Cells.Item(196, 6).Value = "depressed"
$ws.Cells.Item(199, 6).Value = "depressed"
$ws.Cells.Item(200, 6).Value = "depressed"
$ws.Cells.Item(201, 6).Value = "displeased"
